# Applies:
#  - Summary sheet: set column widths (A:21 B:14 C:21 D:19 E:20),
#    rename D1/E1 headers to include "(CZK)"
#  - Czech Data / Poland Data sheets: set column widths
#    (A:21 B:21 C:21 D:14 E:78/73 F:21) and delete the trailing
#    SUM-totals row (row 5)

$wb = $excel.ActiveWorkbook

# Excel's ColumnWidth property is expressed in "characters" of the Normal
# style font and gets re-quantized to pixels on write, which shifts an
# integer input by +5/6 in the saved <col width="..."> value. Subtracting
# 5/6 before assigning keeps the persisted width an exact integer.
function Set-ExactColumnWidth($ws, $colIndex, $targetWidth) {
    $ws.Columns.Item($colIndex).ColumnWidth = $targetWidth - (5 / 6)
}

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")

Set-ExactColumnWidth $wsSummary 1 21
Set-ExactColumnWidth $wsSummary 2 14
Set-ExactColumnWidth $wsSummary 3 21
Set-ExactColumnWidth $wsSummary 4 19
Set-ExactColumnWidth $wsSummary 5 20

$wsSummary.Range("D1").Value = "Czech Price (CZK)"
$wsSummary.Range("E1").Value = "Poland Price (CZK)"

# --- Czech Data sheet ---
$wsCzech = $wb.Worksheets.Item("Czech Data")

Set-ExactColumnWidth $wsCzech 1 21
Set-ExactColumnWidth $wsCzech 2 21
Set-ExactColumnWidth $wsCzech 3 21
Set-ExactColumnWidth $wsCzech 4 14
Set-ExactColumnWidth $wsCzech 5 78
Set-ExactColumnWidth $wsCzech 6 21

$wsCzech.Rows.Item(5).Delete()

# --- Poland Data sheet ---
$wsPoland = $wb.Worksheets.Item("Poland Data")

Set-ExactColumnWidth $wsPoland 1 21
Set-ExactColumnWidth $wsPoland 2 21
Set-ExactColumnWidth $wsPoland 3 21
Set-ExactColumnWidth $wsPoland 4 14
Set-ExactColumnWidth $wsPoland 5 73
Set-ExactColumnWidth $wsPoland 6 21

$wsPoland.Rows.Item(5).Delete()

Write-Output "Applied column widths, header renames and totals-row removal."
